$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the front: old A (Code) -> B, old B (Description) -> C,
# old C (Definition, empty) -> D.
$ws.Columns.Item(1).Insert()

# New column A is the "Version" column.
$ws.Range("A1").Value = "Version"

# Fill the Version value ("1.0") for every data row (2-22).
# Entering it with a leading apostrophe forces it to be stored as text
# (otherwise "1.0" would be parsed as the number 1), then resetting the
# cell style back to Normal drops the quote-prefix formatting so the cell
# ends up as a plain shared-string cell with no style override.
for ($row = 2; $row -le 22; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = "'1.0"
    $cell.Style = "Normal"
}
